$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule values (rows 2-5), plus a new row 6
$data = @(
    @(1, 2, 9, 3, 4, 1, -5, 21),
    @(2, 0, 9, 2, 5, 2, -4, 32),
    @(3, 4, 5, 9, 4, 5, -1, 65),
    @(4, 3, 8, 6, 5, 3, -3, 43),
    @(5, 1, 6, 5, 4, 4, -2, 54)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $ws.Cells.Item($rowIndex, 9).Value = 5
    $ws.Cells.Item($rowIndex, 10).Value = "train_dim2_1"
    $rowIndex++
}

$ws.Range("I1").Select()
